# Account Payables Test Plans
# Update the "expected quantity before/after" sample values on the
# POReceipt sheet (columns O/P, row 2) to the new test-plan figures.
# The values must remain text (e.g. "948.0"), matching how the sheet
# already stores its other sample numbers (593.0/603.0/613.0) as text
# via shared strings rather than as numeric literals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("POReceipt")

# QuantityBefore sample (O2): 603.0 -> 948.0
$ws.Range("O2").NumberFormat = "@"
$ws.Range("O2").Value = "948.0"
$ws.Range("O2").ClearFormats()

# QuantityAfter sample (P2): 613.0 -> 958.0
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "958.0"
$ws.Range("P2").ClearFormats()
